$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-converted to numbers by Excel
# are temporarily formatted as Text ("@") before assignment, then the format
# is cleared again so the cell keeps the workbook default (no explicit style).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '307.45'
$ws.Range("D5").ClearFormats()
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9973'
$ws.Range("D6").ClearFormats()
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3619'
$ws.Range("D7").ClearFormats()
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '47.47'
$ws.Range("D8").ClearFormats()
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3283'
$ws.Range("D9").ClearFormats()
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.128'
$ws.Range("D10").ClearFormats()
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.06953'
$ws.Range("D11").ClearFormats()
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.9965'
$ws.Range("D12").ClearFormats()
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.955'
$ws.Range("D13").ClearFormats()
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '19.30'
$ws.Range("D14").ClearFormats()
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.626'
$ws.Range("D15").ClearFormats()
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001041'
$ws.Range("D17").ClearFormats()
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06507'
$ws.Range("D18").ClearFormats()
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.9980'
$ws.Range("D19").ClearFormats()
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '76.67'
$ws.Range("D20").ClearFormats()
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.931'
$ws.Range("D21").ClearFormats()
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '15.74'
$ws.Range("D22").ClearFormats()
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.68'
$ws.Range("D23").ClearFormats()
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.422'
$ws.Range("D25").ClearFormats()
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.352'
$ws.Range("D26").ClearFormats()
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '146.41'
$ws.Range("D27").ClearFormats()
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.39'
$ws.Range("D28").ClearFormats()
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '124.38'
$ws.Range("D30").ClearFormats()
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.185'
$ws.Range("D31").ClearFormats()
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.039'
$ws.Range("D32").ClearFormats()
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.638'
$ws.Range("D33").ClearFormats()
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08329'
$ws.Range("D34").ClearFormats()
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.667'
$ws.Range("D35").ClearFormats()
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '12.30'
$ws.Range("D36").ClearFormats()
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.238'
$ws.Range("D37").ClearFormats()
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.06040'
$ws.Range("D38").ClearFormats()
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02208'
$ws.Range("D39").ClearFormats()
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.209'
$ws.Range("D40").ClearFormats()
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.2052'
$ws.Range("D41").ClearFormats()
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.183'
$ws.Range("D42").ClearFormats()
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9976'
$ws.Range("D43").ClearFormats()
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5845'
$ws.Range("D44").ClearFormats()
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.733'
$ws.Range("D45").ClearFormats()
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.56'
$ws.Range("D46").ClearFormats()
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5575'
$ws.Range("D47").ClearFormats()
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '122.07'
$ws.Range("D48").ClearFormats()
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.942'
$ws.Range("D49").ClearFormats()
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06911'
$ws.Range("D50").ClearFormats()
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '74.45'
$ws.Range("D51").ClearFormats()

# Remaining cells are safe to assign directly (they are not number-like strings).
$ws.Range("D2").Value = '24.416.31'
$ws.Range("E2").Value = '  -1.93%  '
$ws.Range("D3").Value = '1.653.39'
$ws.Range("E3").Value = '  -3.64%  '
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("E5").Value = '  -1.40%  '
$ws.Range("E6").Value = '  -0.59%  '
$ws.Range("E7").Value = '  -4.18%  '
$ws.Range("E8").Value = '  -4.42%  '
$ws.Range("E9").Value = '  -5.88%  '
$ws.Range("E10").Value = '  -5.56%  '
$ws.Range("E11").Value = '  -7.18%  '
$ws.Range("E12").Value = '  -1.32%  '
$ws.Range("E13").Value = '  -5.40%  '
$ws.Range("E14").Value = '  -7.66%  '
$ws.Range("E15").Value = '  -5.31%  '
$ws.Range("D16").Value = '1.656.64'
$ws.Range("E16").Value = '  -3.70%  '
$ws.Range("E17").Value = '  -7.61%  '
$ws.Range("E18").Value = '  -3.48%  '
$ws.Range("E19").Value = '  -0.50%  '
$ws.Range("E20").Value = '  -9.01%  '
$ws.Range("E21").Value = '  -7.39%  '
$ws.Range("E22").Value = '  -8.71%  '
$ws.Range("E23").Value = '  -2.84%  '
$ws.Range("D24").Value = '24.435.53'
$ws.Range("E24").Value = '  -1.65%  '
$ws.Range("E25").Value = '  -0.79%  '
$ws.Range("E26").Value = '  -15.79%  '
$ws.Range("E27").Value = '  -3.09%  '
$ws.Range("E28").Value = '  -10.57%  '
$ws.Range("D29").Value = '1.836.95'
$ws.Range("E29").Value = '  -3.82%  '
$ws.Range("E30").Value = '  -5.86%  '
$ws.Range("E31").Value = '  +0.54%  '
$ws.Range("E32").Value = '  -4.92%  '
$ws.Range("E33").Value = '  -17.19%  '
$ws.Range("E34").Value = '  -5.91%  '
$ws.Range("E35").Value = '  -5.96%  '
$ws.Range("E36").Value = '  -10.89%  '
$ws.Range("E37").Value = '  -6.59%  '
$ws.Range("E38").Value = '  -7.82%  '
$ws.Range("E39").Value = '  -8.22%  '
$ws.Range("E40").Value = '  -5.91%  '
$ws.Range("B41").Value = 'Algorand'
$ws.Range("C41").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("E41").Value = '  -6.81%  '
$ws.Range("B42").Value = 'FraxShare'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("E42").Value = '  -8.92%  '
$ws.Range("E43").Value = '  -0.16%  '
$ws.Range("E44").Value = '  -9.26%  '
$ws.Range("E45").Value = '  -2.65%  '
$ws.Range("E46").Value = '  -10.05%  '
$ws.Range("E47").Value = '  -9.31%  '
$ws.Range("E48").Value = '  -6.18%  '
$ws.Range("E49").Value = '  -9.45%  '
$ws.Range("E50").Value = '  -5.04%  '
$ws.Range("E51").Value = '  -6.79%  '
